$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 346 (existing rows 346..373 shift down to 348..375).
$ws.Rows.Item(346).Insert()
$ws.Rows.Item(346).Insert()

# New row 346
$ws.Range("A346").Value = 10
$ws.Range("B346").Value = "Vega Modelo de Temuco"
$ws.Range("C346").Value = "La Araucanía"
$ws.Range("D346").Value = 45008
$ws.Range("E346").Value = 9
$ws.Range("F346").Value = 100112052
$ws.Range("G346").Value = "Albahaca"
$ws.Range("H346").Value = "Sin especificar"
$ws.Range("I346").Value = "Primera"
$ws.Range("J346").Value = 35
$ws.Range("K346").Value = 5000
$ws.Range("L346").Value = 5000
$ws.Range("M346").Value = 5000
$ws.Range("N346").Value = "$/paquete"
$ws.Range("O346").Value = "Región de La Araucanía"
$ws.Range("P346").Value = 5000
$ws.Range("Q346").Value = 1
$ws.Range("R346").Value = "Hortaliza"

# New row 347
$ws.Range("A347").Value = 10
$ws.Range("B347").Value = "Vega Modelo de Temuco"
$ws.Range("C347").Value = "La Araucanía"
$ws.Range("D347").Value = 45008
$ws.Range("E347").Value = 9
$ws.Range("F347").Value = 100112052
$ws.Range("G347").Value = "Albahaca"
$ws.Range("H347").Value = "Sin especificar"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 65
$ws.Range("K347").Value = 6000
$ws.Range("L347").Value = 6000
$ws.Range("M347").Value = 6000
$ws.Range("N347").Value = "$/paquete"
$ws.Range("O347").Value = "Región del Maule"
$ws.Range("P347").Value = 6000
$ws.Range("Q347").Value = 1
$ws.Range("R347").Value = "Hortaliza"
